$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 8649
$ws.Range("J69").Value = 8649
$ws.Range("L69").Value = 25947
$ws.Range("N69").Value = -27695
$ws.Range("H72").Value = 8649
$ws.Range("J72").Value = 8649
$ws.Range("L72").Value = 77841
$ws.Range("N72").Value = -86577
$ws.Range("H98").Value = 57955.168
$ws.Range("J98").Value = 16106.889
$ws.Range("L98").Value = 16106.889
$ws.Range("N98").Value = -19102.889
$ws.Range("H103").Value = 1669450.6
$ws.Range("J103").Value = 3599.4
$ws.Range("L103").Value = 10798.2
$ws.Range("N103").Value = -11970.2
$ws.Range("H109").Value = 342051260
$ws.Range("J109").Value = 342051260
$ws.Range("L109").Value = 342051260
$ws.Range("N109").Value = -342054034
$ws.Range("H122").Value = 57955.168
$ws.Range("J122").Value = 16106.889
$ws.Range("L122").Value = 48320.667
$ws.Range("N122").Value = -53220.667
$ws.Range("H137").Value = 355676.7
$ws.Range("I137").Value = 439871.34
$ws.Range("K137").Value = 1319614.02
$ws.Range("M137").Value = -1317064.02

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9581.357
$ws.Range("I2").Value = 10988.728
$ws.Range("K2").Value = 10988.728
$ws.Range("M2").Value = -10875.728
$ws.Range("H32").Value = 7443.515
$ws.Range("I32").Value = 7698.0967
$ws.Range("K32").Value = 7698.0967
$ws.Range("M32").Value = -7411.0967
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").ClearContents()
$ws.Range("N76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").ClearContents()
$ws.Range("N79").Value = 0
$ws.Range("H102").Value = 4509.5386
$ws.Range("I102").Value = 2708.2222
$ws.Range("K102").Value = 2708.2222
$ws.Range("M102").Value = -1086.2222
$ws.Range("H116").Value = 9581.357
$ws.Range("I116").Value = 10988.728
$ws.Range("K116").Value = 10988.728
$ws.Range("M116").Value = -8694.727999999999
$ws.Range("H122").Value = 612332.2
$ws.Range("I122").Value = 2567.9473
$ws.Range("K122").Value = 7703.841899999999
$ws.Range("M122").Value = -5253.841899999999
$ws.Range("H132").Value = 3666.7222
$ws.Range("I132").Value = 2028.3334
$ws.Range("K132").Value = 6085.0002
$ws.Range("M132").Value = -3555.0002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9581.357
$ws.Range("I3").Value = 10988.728
$ws.Range("K3").Value = 10988.728
$ws.Range("M3").Value = -10874.728
$ws.Range("H18").Value = 11000
$ws.Range("J18").Value = 11000
$ws.Range("L18").Value = 11000
$ws.Range("N18").Value = -12058

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4651.5713
$ws.Range("I31").Value = 3485.2354
$ws.Range("J31").Value = 6454.091
$ws.Range("K31").Value = 3485.2354
$ws.Range("L31").Value = 6454.091
$ws.Range("M31").Value = -3190.2354
$ws.Range("N31").Value = -7044.091
$ws.Range("H34").Value = 4651.5713
$ws.Range("I34").Value = 3485.2354
$ws.Range("J34").Value = 6454.091
$ws.Range("K34").Value = 3485.2354
$ws.Range("L34").Value = 6454.091
$ws.Range("M34").Value = -3283.2354
$ws.Range("N34").Value = -6858.091
$ws.Range("H41").Value = 15000
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H43").Value = 21800
$ws.Range("J43").Value = 21800
$ws.Range("L43").Value = 21800
$ws.Range("N43").Value = -22168
$ws.Range("H50").Value = 16500
$ws.Range("H51").Value = 38333.332
$ws.Range("I51").Value = 45000
$ws.Range("K51").Value = 45000
$ws.Range("M51").Value = -44264
$ws.Range("H58").Value = 2984.5667
$ws.Range("I58").Value = 1690.9375
$ws.Range("J58").Value = 4463
$ws.Range("K58").Value = 1690.9375
$ws.Range("L58").Value = 4463
$ws.Range("M58").Value = -1487.9375
$ws.Range("N58").Value = -4869
$ws.Range("H61").Value = 38333.332
$ws.Range("I61").Value = 45000
$ws.Range("K61").Value = 45000
$ws.Range("M61").Value = -44652
$ws.Range("H68").Value = 54749.5
$ws.Range("J68").Value = 54749.5
$ws.Range("L68").Value = 54749.5
$ws.Range("N68").Value = -56247.5
$ws.Range("H71").Value = 54749.5
$ws.Range("J71").Value = 54749.5
$ws.Range("L71").Value = 164248.5
$ws.Range("N71").Value = -171736.5
$ws.Range("H95").Value = 208030800
$ws.Range("J95").Value = 208030800
$ws.Range("L95").Value = 208030800
$ws.Range("N95").Value = -208036292
$ws.Range("H96").Value = 208030400
$ws.Range("J96").Value = 208030400
$ws.Range("L96").Value = 208030400
$ws.Range("N96").Value = -208035892
$ws.Range("H101").Value = 21800
$ws.Range("J101").Value = 21800
$ws.Range("L101").Value = 21800
$ws.Range("N101").Value = -28290
$ws.Range("H102").Value = 77494.5
$ws.Range("I102").Value = 74990
$ws.Range("J102").Value = 79999
$ws.Range("K102").Value = 74990
$ws.Range("L102").Value = 79999
$ws.Range("M102").Value = -72556
$ws.Range("N102").Value = -84867
$ws.Range("H103").Value = 99999
$ws.Range("I103").Value = 99999
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 99999
$ws.Range("L103").Value = 0
$ws.Range("M103").ClearContents()
$ws.Range("N103").Value = -98827
$ws.Range("H107").Value = 10838.375
$ws.Range("J107").Value = 2114.5715
$ws.Range("L107").Value = 2114.5715
$ws.Range("N107").Value = -5954.5715
$ws.Range("H108").Value = 23745.666
$ws.Range("H133").Value = 57000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 57000
$ws.Range("K133").Value = 0
$ws.Range("L133").ClearContents()
$ws.Range("M133").Value = 57000
$ws.Range("N133").Value = -62060
$ws.Range("H136").Value = 2984.5667
$ws.Range("I136").Value = 1690.9375
$ws.Range("J136").Value = 4463
$ws.Range("K136").Value = 5072.8125
$ws.Range("L136").Value = 13389
$ws.Range("M136").Value = -2522.8125
$ws.Range("N136").Value = -18489

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7030.375
$ws.Range("I56").Value = 7030.375
$ws.Range("K56").Value = 7030.375
$ws.Range("M56").Value = -6500.375
$ws.Range("H136").Value = 2893.75
$ws.Range("I136").Value = 2893.75
$ws.Range("K136").Value = 8681.25
$ws.Range("M136").Value = -3581.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4439
$ws.Range("I70").Value = 4233.5
$ws.Range("J70").Value = 4781.5
$ws.Range("K70").Value = 4233.5
$ws.Range("L70").Value = 4781.5
$ws.Range("M70").Value = -3963.5
$ws.Range("N70").Value = -5321.5
$ws.Range("H73").Value = 4439
$ws.Range("I73").Value = 4233.5
$ws.Range("J73").Value = 4781.5
$ws.Range("K73").Value = 4233.5
$ws.Range("L73").Value = 4781.5
$ws.Range("M73").Value = -3297.5
$ws.Range("N73").Value = -6653.5
$ws.Range("H113").Value = 2625.1052
$ws.Range("J113").Value = 2774.8572
$ws.Range("L113").Value = 2774.8572
$ws.Range("N113").Value = -7114.8572
$ws.Range("H141").Value = 100429
$ws.Range("J141").Value = 100429
$ws.Range("L141").Value = 100429
$ws.Range("N141").Value = -110789

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9392862
$ws.Range("I2").Value = 4178580.8
$ws.Range("J2").Value = 15650000
$ws.Range("K2").Value = 4178580.8
$ws.Range("L2").Value = 15650000
$ws.Range("M2").Value = -4178468.8
$ws.Range("N2").Value = -15650224
$ws.Range("H82").Value = 3532.0833
$ws.Range("I82").Value = 4123.7144
$ws.Range("K82").Value = 4123.7144
$ws.Range("M82").Value = -3762.7144
$ws.Range("H85").Value = 3532.0833
$ws.Range("I85").Value = 4123.7144
$ws.Range("K85").Value = 4123.7144
$ws.Range("M85").Value = -2875.7144

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10820000
$ws.Range("I5").Value = 10820000
$ws.Range("K5").Value = 10820000
$ws.Range("M5").Value = -10819888
$ws.Range("H75").Value = 25000
$ws.Range("J75").Value = 25000
$ws.Range("L75").Value = 25000
$ws.Range("N75").Value = -26872
$ws.Range("H78").Value = 25000
$ws.Range("J78").Value = 25000
$ws.Range("L78").Value = 75000
$ws.Range("N78").Value = -84360
$ws.Range("H107").Value = 30684.637
$ws.Range("I107").Value = 3355.2856
$ws.Range("K107").Value = 10065.8568
$ws.Range("M107").Value = -8145.856800000001
$ws.Range("H126").Value = 21301.227
$ws.Range("I126").Value = 25096.295
$ws.Range("K126").Value = 75288.88499999999
$ws.Range("M126").Value = -72818.88499999999
